# Fruta / hortaliza, semanal
# Insert two new weekly records at the top of the data block (rows 39-40),
# pushing the existing rows 39-60 down to rows 41-62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("39:40").Insert()

# --- Row 39: Lapins, Primera ---
$ws.Cells.Item(39, 1).Value = 4
$ws.Cells.Item(39, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(39, 3).Value = "Los Lagos"
$ws.Cells.Item(39, 4).Value = 44572
$ws.Cells.Item(39, 5).Value = 10
$ws.Cells.Item(39, 6).Value = "Fruta"
$ws.Cells.Item(39, 7).Value = 100103
$ws.Cells.Item(39, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(39, 9).Value = 100103001
$ws.Cells.Item(39, 10).Value = "Cereza"
$ws.Cells.Item(39, 11).Value = "Lapins"
$ws.Cells.Item(39, 12).Value = "Primera"
$ws.Cells.Item(39, 13).Value = 600
$ws.Cells.Item(39, 14).Value = 8000
$ws.Cells.Item(39, 15).Value = 8500
$ws.Cells.Item(39, 16).Value = 8250
$ws.Cells.Item(39, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(39, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(39, 19).Value = 825
$ws.Cells.Item(39, 20).Value = 10

# --- Row 40: Santina, Primera ---
$ws.Cells.Item(40, 1).Value = 4
$ws.Cells.Item(40, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(40, 3).Value = "Los Lagos"
$ws.Cells.Item(40, 4).Value = 44572
$ws.Cells.Item(40, 5).Value = 10
$ws.Cells.Item(40, 6).Value = "Fruta"
$ws.Cells.Item(40, 7).Value = 100103
$ws.Cells.Item(40, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(40, 9).Value = 100103001
$ws.Cells.Item(40, 10).Value = "Cereza"
$ws.Cells.Item(40, 11).Value = "Santina"
$ws.Cells.Item(40, 12).Value = "Primera"
$ws.Cells.Item(40, 13).Value = 600
$ws.Cells.Item(40, 14).Value = 7500
$ws.Cells.Item(40, 15).Value = 8000
$ws.Cells.Item(40, 16).Value = 7750
$ws.Cells.Item(40, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(40, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(40, 19).Value = 775
$ws.Cells.Item(40, 20).Value = 10
